# ran resolve and classify+summarise steps after changes to mapping file
$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet: Range Analysis no longer computed -> zero out
#     species counts and drop the now-meaningless percentage column (C).
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()
$wsRange.Range("C3").ClearContents()
$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()
$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()
$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()
$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- "Species qualification" sheet: Range Analysis species count -> 0
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("B5").Value = 0

# --- "High Priority break-up" sheet: re-summarised percentages
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsBreakup.Range("E2").Value = 9.1
$wsBreakup.Range("D3").Value = 10
$wsBreakup.Range("E3").Value = 90.90000000000001
